# "Transpose columns": insert a new "Description" column right after
# ItemNo (A), pushing the existing date/quantity columns from B:G to C:H.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B:B").Insert()

# Populate the new Description column with item descriptions.
$ws.Range("B1").Value2 = "Description"
$ws.Range("B2").Value2 = "This is item 1"
$ws.Range("B3").Value2 = "This is item 2"

# Give the new column the same look as the ItemNo column (A), and a
# width sized for its longer text.
$ws.Range("A1:A3").Copy()
$ws.Range("B1:B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Columns("B:B").ColumnWidth = 12.666666666666666

Write-Host "Done"
